$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 83.5
$ws.Range("C3").Value = 86
$ws.Range("C4").Value = 88
$ws.Range("C6").Value = 85
$ws.Range("C13").Value = 84.59999999999999
$ws.Range("C14").Value = 85
$ws.Range("C17").Value = 83.5
$ws.Range("C18").Value = 88.2
